$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in this week's test-run results (row 5): date + passed/failed/blocked counts.
# Column E already holds a shared SUM formula (E5:E22), so it recalculates itself.
$ws.Range("A5").Value = 41661   # 23-Jan-2018 (date1904 serial)
$ws.Range("B5").Value = 9       # passed
$ws.Range("C5").Value = 1       # failed
$ws.Range("D5").Value = 5       # blocked

# Move the active cell/selection to D5, matching where the author left off editing.
$ws.Range("D5").Select()
